$wb = $excel.ActiveWorkbook

# --- Rename the "Requested quantity" headers ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "PO Forecast"

# Header row
$newSheet.Cells.Item(1,1).Value = "ds"
$newSheet.Cells.Item(1,2).Value = "PO_Forecast"
$newSheet.Cells.Item(1,3).Value = "yhat_lower"
$newSheet.Cells.Item(1,4).Value = "yhat_upper"

# Match the header style (bold / centered / bordered) used on the other sheets
$wsWeekly.Range("B1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)

# Forecast data rows
$newSheet.Cells.Item(2,1).Value = 45417.99999999999
$newSheet.Cells.Item(2,2).Value = 62
$newSheet.Cells.Item(2,3).Value = -8.341637082111493
$newSheet.Cells.Item(2,4).Value = 135.4841267117727
$newSheet.Cells.Item(3,1).Value = 45424.99999999999
$newSheet.Cells.Item(3,2).Value = 62
$newSheet.Cells.Item(3,3).Value = -0.8304351850193069
$newSheet.Cells.Item(3,4).Value = 132.4019540677261
$newSheet.Cells.Item(4,1).Value = 45431.99999999999
$newSheet.Cells.Item(4,2).Value = 62
$newSheet.Cells.Item(4,3).Value = -15.00133043856562
$newSheet.Cells.Item(4,4).Value = 127.3094742931088
$newSheet.Cells.Item(5,1).Value = 45438.99999999999
$newSheet.Cells.Item(5,2).Value = 62
$newSheet.Cells.Item(5,3).Value = -6.510996746614118
$newSheet.Cells.Item(5,4).Value = 134.4728684420786
$newSheet.Cells.Item(6,1).Value = 45445.99999999999
$newSheet.Cells.Item(6,2).Value = 62
$newSheet.Cells.Item(6,3).Value = -17.63348903371085
$newSheet.Cells.Item(6,4).Value = 128.6517333763463
$newSheet.Cells.Item(7,1).Value = 45459.99999999999
$newSheet.Cells.Item(7,2).Value = 63
$newSheet.Cells.Item(7,3).Value = -12.58900209341426
$newSheet.Cells.Item(7,4).Value = 137.4240672501215
$newSheet.Cells.Item(8,1).Value = 45466.99999999999
$newSheet.Cells.Item(8,2).Value = 63
$newSheet.Cells.Item(8,3).Value = -5.672657188729431
$newSheet.Cells.Item(8,4).Value = 137.2905032670542
$newSheet.Cells.Item(9,1).Value = 45473.99999999999
$newSheet.Cells.Item(9,2).Value = 63
$newSheet.Cells.Item(9,3).Value = -12.79626079137827
$newSheet.Cells.Item(9,4).Value = 139.7842170726684
$newSheet.Cells.Item(10,1).Value = 45480.99999999999
$newSheet.Cells.Item(10,2).Value = 63
$newSheet.Cells.Item(10,3).Value = -8.685313729763386
$newSheet.Cells.Item(10,4).Value = 135.3683407176881
$newSheet.Cells.Item(11,1).Value = 45501.99999999999
$newSheet.Cells.Item(11,2).Value = 63
$newSheet.Cells.Item(11,3).Value = -6.790896721986409
$newSheet.Cells.Item(11,4).Value = 139.6006346702594
$newSheet.Cells.Item(12,1).Value = 45508.99999999999
$newSheet.Cells.Item(12,2).Value = 63
$newSheet.Cells.Item(12,3).Value = -12.44435532869558
$newSheet.Cells.Item(12,4).Value = 139.6391033018615
$newSheet.Cells.Item(13,1).Value = 45522.99999999999
$newSheet.Cells.Item(13,2).Value = 63
$newSheet.Cells.Item(13,3).Value = -6.729412347371183
$newSheet.Cells.Item(13,4).Value = 137.8360860151657
$newSheet.Cells.Item(14,1).Value = 45529.99999999999
$newSheet.Cells.Item(14,2).Value = 64
$newSheet.Cells.Item(14,3).Value = -6.825246378833727
$newSheet.Cells.Item(14,4).Value = 133.9884401189239
$newSheet.Cells.Item(15,1).Value = 45536.99999999999
$newSheet.Cells.Item(15,2).Value = 64
$newSheet.Cells.Item(15,3).Value = -8.359587563616961
$newSheet.Cells.Item(15,4).Value = 132.493989553559
$newSheet.Cells.Item(16,1).Value = 45550.99999999999
$newSheet.Cells.Item(16,2).Value = 64
$newSheet.Cells.Item(16,3).Value = -3.918994873437521
$newSheet.Cells.Item(16,4).Value = 139.8139987782872
$newSheet.Cells.Item(17,1).Value = 45564.99999999999
$newSheet.Cells.Item(17,2).Value = 64
$newSheet.Cells.Item(17,3).Value = -7.056710353310593
$newSheet.Cells.Item(17,4).Value = 133.4457399466068
$newSheet.Cells.Item(18,1).Value = 45571.99999999999
$newSheet.Cells.Item(18,2).Value = 64
$newSheet.Cells.Item(18,3).Value = -14.89571647194829
$newSheet.Cells.Item(18,4).Value = 138.8290323152997
$newSheet.Cells.Item(19,1).Value = 45578.99999999999
$newSheet.Cells.Item(19,2).Value = 64
$newSheet.Cells.Item(19,3).Value = -9.808263629709598
$newSheet.Cells.Item(19,4).Value = 132.687095306663
$newSheet.Cells.Item(20,1).Value = 45592.99999999999
$newSheet.Cells.Item(20,2).Value = 64
$newSheet.Cells.Item(20,3).Value = -5.035023971857426
$newSheet.Cells.Item(20,4).Value = 138.5785588108151
$newSheet.Cells.Item(21,1).Value = 45599.99999999999
$newSheet.Cells.Item(21,2).Value = 65
$newSheet.Cells.Item(21,3).Value = -3.52576268175773
$newSheet.Cells.Item(21,4).Value = 133.832210453073
$newSheet.Cells.Item(22,1).Value = 45606.99999999999
$newSheet.Cells.Item(22,2).Value = 65
$newSheet.Cells.Item(22,3).Value = -4.097794756980619
$newSheet.Cells.Item(22,4).Value = 132.3042869642064
$newSheet.Cells.Item(23,1).Value = 45613.99999999999
$newSheet.Cells.Item(23,2).Value = 65
$newSheet.Cells.Item(23,3).Value = -9.990523638419873
$newSheet.Cells.Item(23,4).Value = 129.5946483890869
$newSheet.Cells.Item(24,1).Value = 45620.99999999999
$newSheet.Cells.Item(24,2).Value = 65
$newSheet.Cells.Item(24,3).Value = -6.103859114664812
$newSheet.Cells.Item(24,4).Value = 139.5449085667339
$newSheet.Cells.Item(25,1).Value = 45627.99999999999
$newSheet.Cells.Item(25,2).Value = 65
$newSheet.Cells.Item(25,3).Value = -2.495082545260319
$newSheet.Cells.Item(25,4).Value = 134.4833189082588
$newSheet.Cells.Item(26,1).Value = 45634.99999999999
$newSheet.Cells.Item(26,2).Value = 65
$newSheet.Cells.Item(26,3).Value = -5.697886615359137
$newSheet.Cells.Item(26,4).Value = 135.871284730655
$newSheet.Cells.Item(27,1).Value = 45641.99999999999
$newSheet.Cells.Item(27,2).Value = 65
$newSheet.Cells.Item(27,3).Value = -3.414728552328013
$newSheet.Cells.Item(27,4).Value = 135.0654038454653
$newSheet.Cells.Item(28,1).Value = 45648.99999999999
$newSheet.Cells.Item(28,2).Value = 65
$newSheet.Cells.Item(28,3).Value = -0.621338799852109
$newSheet.Cells.Item(28,4).Value = 137.4681044102001

# Match the date-column style used on the other sheets
$wsWeekly.Range("A2").Copy()
$newSheet.Range("A2:A28").PasteSpecial(-4122)
